$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "29.070.83"
$ws.Range("E2").Value = "  -1.39%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "1.830.07"
$ws.Range("E3").Value = "  -1.41%  "

# Row 4: TetherUSD
$ws.Range("D4").Value = "'0.9989"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.06%  "

# Row 5: BNB
$ws.Range("D5").Value = "'239.01"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.50%  "

# Row 6: XRP
$ws.Range("D6").Value = "'0.6610"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -4.75%  "

# Row 7: USDC
$ws.Range("D7").Value = "'0.9996"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.04%  "

# Row 8: Cardano
$ws.Range("D8").Value = "'0.2942"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -3.91%  "

# Row 9: Dogecoin
$ws.Range("D9").Value = "'0.07327"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -4.37%  "

# Row 10: Solana
$ws.Range("D10").Value = "'22.65"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -3.97%  "

# Row 11: TRON
$ws.Range("D11").Value = "'0.07640"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.69%  "

# Row 12: WrappedEther
$ws.Range("D12").Value = "1.841.09"
$ws.Range("E12").Value = "  -0.65%  "

# Row 13: Polkadot
$ws.Range("D13").Value = "'5.004"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -2.57%  "

# Row 14: Polygon
$ws.Range("D14").Value = "'0.6720"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.69%  "

# Row 15: Litecoin
$ws.Range("D15").Value = "'86.04"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -5.31%  "

# Row 16: Uniswap
$ws.Range("D16").Value = "'6.101"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -3.93%  "

# Row 17: WrappedBTC
$ws.Range("D17").Value = "29.064.49"
$ws.Range("E17").Value = "  -1.40%  "

# Row 18: ShibaInu
$ws.Range("D18").Value = "'0.000008179"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.11%  "

# Row 19: BitcoinCash
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").Value = "'227.60"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -4.23%  "

# Row 20: Avalanche
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").Value = "'12.45"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.12%  "

# Row 21: Dai
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value = "'0.9995"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.03%  "

# Row 22: Chainlink
$ws.Range("B22").Value = "Chainlink"
$ws.Range("C22").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D22").Value = "'7.239"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -5.37%  "

# Row 23: BinanceUSD
$ws.Range("B23").Value = "BinanceUSD"
$ws.Range("C23").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D23").Value = "'0.9997"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.05%  "

# Row 24: Monero
$ws.Range("B24").Value = "Monero"
$ws.Range("C24").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D24").Value = "'160.66"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.65%  "

# Row 25: Stellar
$ws.Range("B25").Value = "Stellar"
$ws.Range("C25").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D25").Value = "'0.1419"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -4.85%  "

# Row 26: Cosmos
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Value = "'8.633"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.94%  "

# Row 27: EthereumClassic
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'17.92"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.84%  "

# Row 28: PancakeSwap
$ws.Range("B28").Value = "PancakeSwap"
$ws.Range("C28").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D28").Value = "'1.496"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.60%  "

# Row 29: Filecoin
$ws.Range("B29").Value = "Filecoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D29").Value = "'4.213"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.71%  "

# Row 30: InternetComputer(DFINITY)
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").Value = "'4.096"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.44%  "

# Row 31: Toncoin
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").Value = "'1.197"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.58%  "

# Row 32: Hedera
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Value = "'0.05295"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +3.78%  "

# Row 33: ImmutableX
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value = "'0.7478"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -3.10%  "

# Row 34: LidoDAOToken
$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").Value = "'1.846"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.33%  "

# Row 35: ARBITRUM
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "'1.125"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -2.12%  "

# Row 36: HuobiToken
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "'2.681"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.09%  "

# Row 37: Maker
$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D37").Value = "1.294.16"
$ws.Range("E37").Value = "  -2.97%  "

# Row 38: VeChain
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.01803"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.34%  "

# Row 39: MXToken
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "'2.703"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.45%  "

# Row 40: TrustWalletToken
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "'0.9205"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.82%  "

# Row 41: FraxShare
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "'5.965"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +2.14%  "

# Row 42: PaxDollar
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").Value = "'0.9985"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.17%  "

# Row 43: Quant
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "'103.46"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.51%  "

# Row 44: RocketPoolETH
$ws.Range("B44").Value = "RocketPoolETH"
$ws.Range("C44").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D44").Value = "1.984.16"
$ws.Range("E44").Value = "  -0.79%  "

# Row 45: BabyDogeCoin
$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").Value = "'0.00000000123"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.00%  "

# Row 46: Mantle
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").Value = "'0.5174"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.87%  "

# Row 47: RenderToken
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").Value = "'1.746"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.95%  "

# Row 48: Aave
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "'63.14"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.05%  "

# Row 49: EnergySwap
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'9.220"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -5.90%  "

# Row 50: XinFinNetwork
$ws.Range("D50").Value = "'0.07494"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +11.88%  "

# Row 51: Cronos
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.05907"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.27%  "

